$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Q62:Q65 flip from 1 to 0 ---
$ws.Cells.Item(62, 17).Value = 0
$ws.Cells.Item(63, 17).Value = 0
$ws.Cells.Item(64, 17).Value = 0
$ws.Cells.Item(65, 17).Value = 0

# --- 2) O1146 changes from 0 to 2 ---
$ws.Cells.Item(1146, 15).Value = 2

# --- 3) R1148 and R1149 change from empty inline string to numeric 0 ---
$ws.Cells.Item(1148, 18).Value = 0
$ws.Cells.Item(1149, 18).Value = 0

# --- 4) Append new weekly rows 1150-1166 ---
# Columns: A Datetime, B Open, C High, D Low, E Close, F Adj Close, G Volume,
#          H Year, I Month, J Day, K Hour, L Minute, M Second, N Week,
#          O isPivot, P two_line_structure, Q detect_structure
$newRows = @(
    @(45474, 504.2000122070312,  519.5,               503.7000122070312,  514.4000244140625,  511.1280517578125,  5718979,  2024, 7,  1, 0, 0, 0, 27, 0, 0, 0),
    @(45481, 514.4000244140625,  532.5,               505.0499877929688,  526,                 522.6542358398438,  11585045, 2024, 7,  8, 0, 0, 0, 28, 0, 0, 0),
    @(45488, 528,                536.8499755859375,   507.2000122070312,  523.75,              520.4185180664062,  8520292,  2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(45495, 523.9000244140625,  543.5,               511.0499877929688,  542.4000244140625,  538.949951171875,   5286433,  2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 545.5,              558.4500122070312,   537.75,             550.25,              546.75,              4904211,  2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(45509, 541.0499877929688,  549.6500244140625,   514.7999877929688,  531.6500244140625,  528.268310546875,   6086206,  2024, 8,  5, 0, 0, 0, 32, 0, 0, 0),
    @(45516, 532.9500122070312,  551.9500122070312,   521,                549.9000244140625,  549.9000244140625,  8447755,  2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(45523, 558.5,              597.2000122070312,   548.0499877929688,  565.7000122070312,  565.7000122070312,  11382511, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 565.7000122070312,  586,                 563.0499877929688,  571.0499877929688,  571.0499877929688,  4259637,  2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(45537, 574.0999755859375,  604.5,               570.1500244140625,  596.2999877929688,  596.2999877929688,  10349165, 2024, 9,  2, 0, 0, 0, 36, 0, 0, 0),
    @(45544, 598,                626.9000244140625,   593.6500244140625,  622.6500244140625,  622.6500244140625,  10531262, 2024, 9,  9, 0, 0, 0, 37, 0, 0, 0),
    @(45551, 621.0999755859375,  629.5,               613.5499877929688,  622.8499755859375,  622.8499755859375,  4677468,  2024, 9, 16, 0, 0, 0, 38, 1, 0, 0),
    @(45558, 623.25,             629.1500244140625,   600.2000122070312,  617.5499877929688,  617.5499877929688,  13830831, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0),
    @(45565, 616.5,              625,                 574.25,             576.9500122070312,  576.9500122070312,  5873654,  2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(45572, 576.9000244140625,  585,                 561.2000122070312,  571.5499877929688,  571.5499877929688,  4376886,  2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(45579, 571,                588.5999755859375,   549.0499877929688,  556.3499755859375,  556.3499755859375,  5101217,  2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(45586, 557.2000122070312,  565.4500122070312,   530.4000244140625,  536.9500122070312,  536.9500122070312,  4011361,  2024, 10, 21, 0, 0, 0, 43, 0, 0, 0)
)

$targetRow = 1150
foreach ($rowValues in $newRows) {
    # Copy the date cell style (s="2") from the row above so the new Datetime
    # cell keeps the same number format as the rest of column A.
    $ws.Cells.Item($targetRow - 1, 1).Copy($ws.Cells.Item($targetRow, 1))

    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $rowValues[$col - 1]
    }
    $targetRow++
}
